$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The reverse-charge lookup rows for rc_type_2 and rc_type_3 had their
# purchase/sale tax references swapped in the refreshed test fixture:
# row 3 (rc_type_2) should now carry the external.a41a* pair, and row 4
# (rc_type_3) should carry the z0bug.tax_a8aa / external.aa8av pair.
$c3 = $ws.Range("C3").Value2
$d3 = $ws.Range("D3").Value2
$c4 = $ws.Range("C4").Value2
$d4 = $ws.Range("D4").Value2

$ws.Range("C3").Value2 = $c4
$ws.Range("D3").Value2 = $d4
$ws.Range("C4").Value2 = $c3
$ws.Range("D4").Value2 = $d3

# The saved view's active cell moved from D2 to B5.
$ws.Range("B5").Select() | Out-Null
